$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 101, shifting existing rows 101:203 down to 102:204.
$ws.Rows("101").Insert()

# Populate the new row 101 with this week's data (matches the style/date format
# already used by the rest of column D - inherited from the row that used to be
# row 101, since Insert() copies formatting from the row above by default in
# Excel; set format explicitly to be safe).
$ws.Range("A101").Value = 1
$ws.Range("B101").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C101").Value = "Arica y Parinacota"
$ws.Range("D101").Value = [DateTime]"2022-03-09"
$ws.Range("D101").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E101").Value = 15
$ws.Range("F101").Value = "Fruta"
$ws.Range("G101").Value = 100102
$ws.Range("H101").Value = "Cítricos"
$ws.Range("I101").Value = 100102003
$ws.Range("J101").Value = "Limón"
$ws.Range("K101").Value = "Sin especificar"
$ws.Range("L101").Value = "2a amarillo"
$ws.Range("M101").Value = 270
$ws.Range("N101").Value = 29000
$ws.Range("O101").Value = 30000
$ws.Range("P101").Value = 29500
$ws.Range("Q101").Value = "`$/caja 20 kilos"
$ws.Range("R101").Value = "Región de Coquimbo"
$ws.Range("S101").Value = 1475
$ws.Range("T101").Value = 20
